# Insert a new price record as row 57 (pushing the existing rows 57..136
# down to 58..137), matching the author's commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(57).Insert()

$ws.Range("A57").Value = 1
$ws.Range("B57").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C57").Value = "Arica y Parinacota"
$ws.Range("D57").Value = 44469
$ws.Range("E57").Value = 15
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100108
$ws.Range("H57").Value = "Tropicales y subtropicales"
$ws.Range("I57").Value = 100108006
$ws.Range("J57").Value = "Plátano"
$ws.Range("K57").Value = "Sin especificar"
$ws.Range("L57").Value = "Pintón"
$ws.Range("M57").Value = 120
$ws.Range("N57").Value = 17000
$ws.Range("O57").Value = 18000
$ws.Range("P57").Value = 17500
$ws.Range("Q57").Value = "$/caja 20 kilos"
$ws.Range("R57").Value = "Ecuador"
$ws.Range("S57").Value = 875
$ws.Range("T57").Value = 20
